# Assessment 1 Team Contribution Sheet - fill in "Final Deliverable" (column G)
# scores for the three team members who already have a complete row, leaving
# the still-incomplete row (row 6 / Josh McQueen, whose Code Review 4 score in
# column F is also still blank) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G4").Value = 33.3
$ws.Range("G5").Value = 33.3
$ws.Range("G7").Value = 33.3

# Reflect the cell that was selected when the sheet was last saved.
$ws.Range("G6").Select()
